$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 61405410
$ws.Range("I62").Value = 89745256
$ws.Range("K62").Value = 89745256
$ws.Range("M62").Value = -89744632
# Row 65
$ws.Range("H65").Value = 61405410
$ws.Range("I65").Value = 89745256
$ws.Range("K65").Value = 448726280
$ws.Range("M65").Value = -448723160
# Row 94
$ws.Range("H94").Value = 5970.25
$ws.Range("I94").Value = 5970.25
$ws.Range("K94").Value = 5970.25
$ws.Range("M94").Value = -5519.25
# Row 131
$ws.Range("H131").Value = 4007.0833
$ws.Range("I131").Value = 4007.7273
$ws.Range("K131").Value = 12023.1819
$ws.Range("M131").Value = -6983.1819
# Row 132
$ws.Range("H132").Value = 1279.3182
$ws.Range("I132").Value = 1283.0952
$ws.Range("K132").Value = 3849.2856
$ws.Range("M132").Value = -1319.2856
# Row 137
$ws.Range("H137").Value = 4551830
$ws.Range("I137").Value = 8624198
$ws.Range("J137").Value = 9573.462
$ws.Range("K137").Value = 25872594
$ws.Range("L137").Value = 28720.386
$ws.Range("M137").Value = -25870044
$ws.Range("N137").Value = -33820.386
# Row 138
$ws.Range("H138").Value = 4235.533
$ws.Range("J138").Value = 4536.85
$ws.Range("L138").Value = 13610.55
$ws.Range("N138").Value = -23890.55

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6287456
$ws.Range("J32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("N32").Value = -20574
# Row 46
$ws.Range("H46").Value = 14205.444
$ws.Range("J46").Value = 21479.4
$ws.Range("L46").Value = 21479.4
$ws.Range("N46").Value = -22117.4
# Row 61
$ws.Range("H61").Value = 3481.2222
$ws.Range("I61").Value = 2083.1428
$ws.Range("K61").Value = 2083.1428
$ws.Range("M61").Value = -1871.1428
# Row 74
$ws.Range("H74").Value = 4926.4287
$ws.Range("I74").Value = 3697.0667
$ws.Range("K74").Value = 3697.0667
$ws.Range("M74").Value = -2823.0667
# Row 77
$ws.Range("H77").Value = 4926.4287
$ws.Range("I77").Value = 3697.0667
$ws.Range("K77").Value = 18485.3335
$ws.Range("M77").Value = -14117.3335
# Row 110
$ws.Range("H110").Value = 12441025
$ws.Range("I110").Value = 18889506
$ws.Range("J110").Value = 4670.7144
$ws.Range("K110").Value = 18889506
$ws.Range("L110").Value = 4670.7144
$ws.Range("M110").Value = -18887461
$ws.Range("N110").Value = -8760.714400000001
# Row 132
$ws.Range("H132").Value = 3793.884
$ws.Range("I132").Value = 2329.0625
$ws.Range("J132").Value = 7142.048
$ws.Range("K132").Value = 6987.1875
$ws.Range("L132").Value = 21426.144
$ws.Range("M132").Value = -4457.1875
$ws.Range("N132").Value = -26486.144
# Row 136
$ws.Range("H136").Value = 3481.2222
$ws.Range("I136").Value = 2083.1428
$ws.Range("K136").Value = 6249.428400000001
$ws.Range("M136").Value = -3699.428400000001

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1296.6154
$ws.Range("I94").Value = 1282.2858
$ws.Range("J94").Value = 1313.3334
$ws.Range("K94").Value = 1282.2858
$ws.Range("L94").Value = 1313.3334
$ws.Range("M94").Value = -831.2858000000001
$ws.Range("N94").Value = -2215.3334
# Row 105
$ws.Range("H105").Value = 19610570
$ws.Range("I105").Value = 27780624
$ws.Range("J105").Value = 2442.3333
$ws.Range("K105").Value = 27780624
$ws.Range("L105").Value = 2442.3333
$ws.Range("M105").Value = -27778877
$ws.Range("N105").Value = -5936.3333
# Row 134
$ws.Range("H134").Value = 4845.6
$ws.Range("I134").Value = 2083.7856
$ws.Range("J134").Value = 8360.637000000001
$ws.Range("K134").Value = 6251.3568
$ws.Range("L134").Value = 25081.911
$ws.Range("M134").Value = -3716.3568
$ws.Range("N134").Value = -30151.911

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 55559640
$ws.Range("I31").Value = 100002160
$ws.Range("J31").Value = 6492.375
$ws.Range("K31").Value = 100002160
$ws.Range("L31").Value = 6492.375
$ws.Range("M31").Value = -100001865
$ws.Range("N31").Value = -7082.375
# Row 34
$ws.Range("H34").Value = 55559640
$ws.Range("I34").Value = 100002160
$ws.Range("J34").Value = 6492.375
$ws.Range("K34").Value = 100002160
$ws.Range("L34").Value = 6492.375
$ws.Range("M34").Value = -100001958
$ws.Range("N34").Value = -6896.375
# Row 132
$ws.Range("H132").Value = 33868.023
$ws.Range("I132").Value = 2634.5144
$ws.Range("K132").Value = 7903.5432
$ws.Range("M132").Value = -5373.5432
# Row 134
$ws.Range("H134").Value = 5236.875
$ws.Range("I134").Value = 4229.8823
$ws.Range("K134").Value = 12689.6469
$ws.Range("M134").Value = -10154.6469

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 128823.44
$ws.Range("I68").Value = 401298.4
$ws.Range("J68").Value = 4971.1816
$ws.Range("K68").Value = 1203895.2
$ws.Range("L68").Value = 14913.5448
$ws.Range("M68").Value = -1203084.2
$ws.Range("N68").Value = -16535.5448
# Row 71
$ws.Range("H71").Value = 128823.44
$ws.Range("I71").Value = 401298.4
$ws.Range("J71").Value = 4971.1816
$ws.Range("K71").Value = 3611685.6
$ws.Range("L71").Value = 44740.6344
$ws.Range("M71").Value = -3607629.6
$ws.Range("N71").Value = -52852.6344
# Row 80
$ws.Range("H80").Value = 3966.1667
$ws.Range("J80").Value = 3966.1667
$ws.Range("L80").Value = 11898.5001
$ws.Range("N80").Value = -13770.5001
# Row 83
$ws.Range("H83").Value = 3966.1667
$ws.Range("J83").Value = 3966.1667
$ws.Range("L83").Value = 35695.5003
$ws.Range("N83").Value = -45055.5003
# Row 88
$ws.Range("H88").Value = 1900
$ws.Range("I88").Value = 1800
$ws.Range("K88").Value = 5400
$ws.Range("M88").Value = -4972
# Row 91
$ws.Range("H91").Value = 1900
$ws.Range("I91").Value = 1800
$ws.Range("K91").Value = 5400
$ws.Range("M91").Value = -3918

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 851.75757
$ws.Range("I97").Value = 875.75
$ws.Range("J97").Value = 787.7778
$ws.Range("K97").Value = 875.75
$ws.Range("L97").Value = 787.7778
$ws.Range("M97").Value = -379.75
$ws.Range("N97").Value = -1779.7778
# Row 102
$ws.Range("H102").Value = 1324
$ws.Range("I102").Value = 905
$ws.Range("K102").Value = 905
$ws.Range("M102").Value = 717
# Row 107
$ws.Range("H107").Value = 1440.7693
$ws.Range("J107").Value = 1119
$ws.Range("L107").Value = 1119
$ws.Range("N107").Value = -4959
# Row 123
$ws.Range("H123").Value = 54170.184
$ws.Range("J123").Value = 62267.43
$ws.Range("L123").Value = 62267.43
$ws.Range("N123").Value = -67167.42999999999
# Row 126
$ws.Range("H126").Value = 3067.5881
$ws.Range("J126").Value = 5173.2
$ws.Range("L126").Value = 15519.6
$ws.Range("N126").Value = -20459.6
# Row 132
$ws.Range("H132").Value = 3997.0435
$ws.Range("I132").Value = 2640.639
$ws.Range("J132").Value = 8880.1
$ws.Range("K132").Value = 7921.917
$ws.Range("L132").Value = 26640.3
$ws.Range("M132").Value = -5391.917
$ws.Range("N132").Value = -31700.3

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 13162157
$ws.Range("I100").Value = 22730772
$ws.Range("K100").Value = 22730772
$ws.Range("M100").Value = -22730231
# Row 122
$ws.Range("H122").Value = 3379.4
$ws.Range("I122").Value = 2976
$ws.Range("J122").Value = 4993
$ws.Range("K122").Value = 8928
$ws.Range("L122").Value = 14979
$ws.Range("M122").Value = -6478
$ws.Range("N122").Value = -19879
# Row 136
$ws.Range("H136").Value = 5737.615
$ws.Range("I136").Value = 2009.8889
$ws.Range("J136").Value = 7711.1177
$ws.Range("K136").Value = 6029.6667
$ws.Range("L136").Value = 23133.3531
$ws.Range("M136").Value = -3479.6667
$ws.Range("N136").Value = -28233.3531

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 5360.067
$ws.Range("I122").Value = 5275
$ws.Range("J122").Value = 5434.5
$ws.Range("K122").Value = 15825
$ws.Range("L122").Value = 16303.5
$ws.Range("M122").Value = -13375
$ws.Range("N122").Value = -21203.5
# Row 136
$ws.Range("H136").Value = 4853.44
$ws.Range("I136").Value = 2446.8
$ws.Range("K136").Value = 7340.400000000001
$ws.Range("M136").Value = -4790.400000000001
